$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for rows 3 and 4 so we can swap
# their taxon-related data while handling the "Taxonsorteringsordning"
# (column B) specially. Use Value2 (Value has a reflection quirk on this
# host that returns a descriptor string instead of the cell's contents).
$A3 = $ws.Range("A3").Value2
$A4 = $ws.Range("A4").Value2

$E3 = $ws.Range("E3").Value2
$E4 = $ws.Range("E4").Value2

$F3 = $ws.Range("F3").Value2
$F4 = $ws.Range("F4").Value2

$G3 = $ws.Range("G3").Value2
$G4 = $ws.Range("G4").Value2

$H3 = $ws.Range("H3").Value2
$H4 = $ws.Range("H4").Value2

$Q3 = $ws.Range("Q3").Value2
$Q4 = $ws.Range("Q4").Value2

$R3 = $ws.Range("R3").Value2
$R4 = $ws.Range("R4").Value2

$B4old = $ws.Range("B4").Value2

# Column B (Taxonsorteringsordning) for the Leucobryum glaucum taxon was
# updated in the source data from 93539 to 93553.
$newSortOrder = 93553

# B2 keeps the Leucobryum glaucum row, just with the refreshed sort order.
$ws.Range("B2").Value2 = $newSortOrder

# Rows 3 and 4 swap their taxon / location data: what used to be row 4
# (Vanlig groda / Rana temporaria) becomes row 3, and what used to be
# row 3 (Blåmossa / Leucobryum glaucum) becomes row 4.
$ws.Range("A3").Value2 = $A4
$ws.Range("A4").Value2 = $A3

$ws.Range("B3").Value2 = $B4old
$ws.Range("B4").Value2 = $newSortOrder

$ws.Range("E3").Value2 = $E4
$ws.Range("E4").Value2 = $E3

$ws.Range("F3").Value2 = $F4
$ws.Range("F4").Value2 = $F3

$ws.Range("G3").Value2 = $G4
$ws.Range("G4").Value2 = $G3

$ws.Range("H3").Value2 = $H4
$ws.Range("H4").Value2 = $H3

$ws.Range("Q3").Value2 = $Q4
$ws.Range("Q4").Value2 = $Q3

$ws.Range("R3").Value2 = $R4
$ws.Range("R4").Value2 = $R3

# The empty "Bestämningsmetod" marker cell also moves from row 4 to row 3
# (present in row 4, absent in row 3, before the edit — and vice versa
# afterwards): cut it from AF4 into AF3, then re-mark AF3 as an empty text
# cell (rather than a cleared/blank one) and reset its style so it matches
# the original empty-text marker cell exactly.
$ws.Range("AF4").Cut($ws.Range("AF3"))
$ws.Range("AF3").Formula = "'"
$ws.Range("AF3").Style = "Normal"
